$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing headers
$ws.Range("B1").Value = "UnderperformingSKU"
$ws.Range("C1").Value = "UnderperformingMFNPartNumber"

# Update existing comment value
$ws.Range("H2").Value = "test cost structure"

# New headers
$ws.Range("I1").Value = "AgingSKU"
$ws.Range("J1").Value = "AgingMFNPartNumber"
$ws.Range("K1").Value = "UpdateAction"
$ws.Range("L1").Value = "UpdateComment"

# New data values
$ws.Range("I2").Value = "0005352280"
$ws.Range("J2").Value = "DOCK182AUZ"
$ws.Range("K2").Value = "CM - Pricing"
$ws.Range("L2").Value = "test pricing"
